$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021 United Kingdom" header in column O (row 1),
# keeping its existing header style.
$ws.Range("O1").Value = "2021 United Kingdom"

# The new column needs the same green "bingo cell" fill/font formatting
# already used by the other data columns (e.g. column B) in each row
# that has data. Copy that formatting across without touching the values.
$ws.Range("B2").Copy()
$ws.Range("O2:O13").PasteSpecial(-4122)
$ws.Range("O16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
